# Rename the logo pictures that live in the document's headers/footers.
#
#   - The Pearson logo (inserted in both the primary and the first-page
#     footer) is renamed from "image2.png" to "image1.png".
#   - The BTec logo (inserted in the first-page header) is renamed from
#     "image1.jpg" to "image2.jpg".
#
# Word only exposes these pictures through InlineShapes hanging off each
# HeaderFooter's Range - they are not part of the main story, so
# $d.InlineShapes does not see them.

$d = $word.ActiveDocument

foreach ($sec in $d.Sections) {
    # wdHeaderFooterPrimary=1, wdHeaderFooterFirstPage=2, wdHeaderFooterEvenPages=3
    foreach ($idx in 1..3) {
        $ftr = $sec.Footers($idx)
        if ($ftr.Exists) {
            foreach ($shp in $ftr.Range.InlineShapes) {
                if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                    $shp.Name = "image1.png"
                }
            }
        }

        $hdr = $sec.Headers($idx)
        if ($hdr.Exists) {
            foreach ($shp in $hdr.Range.InlineShapes) {
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    $shp.Name = "image2.jpg"
                }
            }
        }
    }
}
